# Insert a new row at position 4 (shifts existing rows 4-17 down to 5-18)
# to add a new fixture entry for image-00002.extra.periods.dcm in the
# Head CT dataset, reflecting the matcher no longer swallowing the
# character after {filename} as a delimiter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

function Set-TextValue {
    param($cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "A4" "Scans"
Set-TextValue "B4" "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Head CT/image-00002.extra.periods.dcm"
Set-TextValue "C4" "image-00002.extra.periods.dcm"
Set-TextValue "D4" "Y"
Set-TextValue "F4" "002304"
Set-TextValue "G4" "20200312"
Set-TextValue "H4" "Head_CT"
Set-TextValue "I4" "002304_CT1"
Set-TextValue "L4" "CT"
$ws.Range("M4").Value = 3
Set-TextValue "O4" "20200312"
Set-TextValue "P4" "CT1 abdomen"
Set-TextValue "Q4" "DOE^JOHN"
Set-TextValue "R4" "002304"
Set-TextValue "S4" "Head CT"
Set-TextValue "T4" "image-00002.extra.periods"

$ws.Range("A4").Select() | Out-Null
